$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 (Joystick): replace the old AliExpress alternate-source link with a
# SparkFun link, and clear the now-stale "cheaper from Aliexpress" note
$ws.Range("G29").Value = "https://www.sparkfun.com/products/9032"

# Row 15 (Teensy 3.2, U1): alternate source now also lists the Adafruit product
# page alongside the existing pjrc.com page
$ws.Range("G15").Value = "https://www.adafruit.com/product/2756 or https://www.pjrc.com/store/teensy32.html"

# Row 9 (C1/C2 electrolytic caps): value "10u" -> "10uF"
$ws.Range("B9").Value = "10uF"

# Clear the stale "cheaper from Aliexpress" note now that the joystick source changed
$ws.Range("H29").Value = ""

# Remove the old "Joystick Cap" row entirely (row 30) - no longer a separate
# line item now that the joystick source note changed
$ws.Rows(30).Delete()

# Update the active selection to match the author's final cursor position
[void]$ws.Range("B10").Select()
